# "update complete, cancel genectics case"
# - Sheet1: the genetics case id (A2) is replaced with a new case id.
# - Sheet2: the status (B6) of the "assign follow up -> approved 3" test
#   case is updated from pass to fail (the case is cancelled).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws1.Range("A2").Value = "CA-FR2PHPWO"

$ws2.Range("B6").Value = "fail"
